$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.192.96'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '1.806.87'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5135'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3956'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07797'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.108'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '40.92'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.376'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.002'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.40'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.809.06'
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.314'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.96'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.91%  '
$ws.Range('E18').Value = '  -2.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06580'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.29'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.011'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = '28.218.71'
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('E24').Value = '  -2.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.214'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.60'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.444'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.51'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('D29').Value = '2.014.09'
$ws.Range('E29').Value = '  -0.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.02'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.91%  '
$ws.Range('E31').Value = '  -0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.059'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.660'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.561'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07150'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.169'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02349'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.037'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.36%  '
$ws.Range('E40').Value = '  -6.28%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6155'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.68%  '
$ws.Range('E42').Value = '  -0.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.151'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.29%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5962'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.67%  '
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.303'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.735'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.82'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.212'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.21%  '
$ws.Range('E50').Value = '  -3.31%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06793'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.38%  '
